$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 1.29
$ws.Range("P2").Value = 1.28
$ws.Range("Q3").Value = 1.71
$ws.Range("Q4").Value = 2.56
$ws.Range("H5").Value = 1.51
$ws.Range("Q9").Value = 1.84
$ws.Range("F10").Value = 3.3
$ws.Range("G10").Value = 5.3
$ws.Range("H10").Value = 1.79
$ws.Range("I10").Value = 2.06
$ws.Range("J10").Value = 3.05
$ws.Range("K10").Value = 5.8
$ws.Range("Q10").Value = 1.61
$ws.Range("P11").Value = 2.52
$ws.Range("Q11").Value = 1.52
$ws.Range("G12").Value = 2.72
$ws.Range("I12").Value = 3.65
$ws.Range("P12").Value = 1.92
$ws.Range("I13").Value = 3.3
$ws.Range("F14").Value = 10
$ws.Range("G14").Value = 14
$ws.Range("H14").Value = 1.34
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 5.7
$ws.Range("F15").Value = 4.1
$ws.Range("G15").Value = 4.9
$ws.Range("H15").Value = 1.72
$ws.Range("I15").Value = 1.88
$ws.Range("J15").Value = 4.2
$ws.Range("K15").Value = 5
$ws.Range("N15").Value = 5.3
$ws.Range("R15").Value = 1.7
$ws.Range("S15").Value = 2.16
$ws.Range("V15").Value = 2.12
$ws.Range("W15").Value = 1.25
$ws.Range("F16").Value = 2.64
$ws.Range("P17").Value = 1.68
$ws.Range("N18").Value = 3.6
$ws.Range("T18").Value = 1.84
$ws.Range("AC18").Value = 7.4
$ws.Range("AD18").Value = 13.5
$ws.Range("AH18").Value = 17
$ws.Range("AJ18").Value = 36
$ws.Range("AK18").Value = 30
$ws.Range("G19").Value = 1.86
$ws.Range("Q20").Value = 1.74
